$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I, shifting existing columns I:N to J:O
$ws.Range("I1").EntireColumn.Insert()

# Populate the new "Onboarding Completed" column
$ws.Range("I1").Value = "Onboarding Completed"
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "Yes"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

$ws.Range("I2").Select()
